# Apply the "LinuxForHealth" rebrand + republish edits described by the
# commit's diff to StructureDefinition-method.xlsx.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Metadata" sheet: StructureDefinition header properties
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/method"

# Version bump: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: new publish timestamp
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------
# "Elements" sheet: the top-level "Extension" row no longer carries the
# ele-1 / ext-1 constraint text in its Constraint(s) column (that text
# stays only on the Extension.extension row below it).
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The Extension.url row's "Fixed Value" cell reuses the same StructureDefinition
# URL string shown on the Metadata sheet, so it picks up the same rebrand.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/method"
